# Convert the text "FALSE" values in the Categories sheet's isMissing column (C)
# into real boolean FALSE values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Categories")

# Replicate the user's selection before performing the edit: the boolean
# column for the data rows (C2:C131), with the view ending on C3:C131.
$ws.Activate()
$ws.Range("C3:C131").Select()

# Replace the text "FALSE" values with actual boolean FALSE values.
$ws.Range("C2:C131").Value = $false
